$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Team Name column (B2:B6) with the new team names
# (order matters for shared-string table insertion order on save)
$ws.Range("B2").Value = "Hackathians"
$ws.Range("B3").Value = "NeuroQuad"
$ws.Range("B5").Value = "RavenClaw"
$ws.Range("B4").Value = "GRYFFINDERS"
$ws.Range("B6").Value = "SLyTHerin"

# Update sheet view: reset scrolled topLeftCell and move selection to E20
$ws.Range("E20").Select()

# Update page setup orientation to portrait
$ws.PageSetup.Orientation = 1
